# Applies the edit described by the upstream commit:
#  1. Three tables (on slides 14, 15 and 16) are switched from the
#     "Table_0" table style to the built-in "No Style, No Grid" table
#     style ({31CDAA57-7DB9-43F2-B19A-7BE244136A39}).
#  2. The presentation's theme (ppt/theme/theme1.xml, the theme used by
#     the slide master) is changed from the "Integral" / "Red Violet"
#     theme to the default Office theme's colour scheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Re-style the three tables that used the old "Table_0" style.
# ---------------------------------------------------------------------
$newTableStyleId = "{31CDAA57-7DB9-43F2-B19A-7BE244136A39}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2. Swap the "Integral" (Red Violet) colour scheme for the default
#    Office theme colours on the slide master's theme.
# ---------------------------------------------------------------------
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388    # dk2      44546A
    4  = 15132391   # lt2      E7E6E6
    5  = 13998939   # accent1  5B9BD5
    6  = 3243501     # accent2  ED7D31
    7  = 10855845   # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308   # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i]
}
